# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 846 for "Vega Monumental Concepción" / Pera,
# shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 846; this shifts rows 846-880 down to 847-881
# and copies the date-format style that row already had in column D.
$ws.Rows.Item(846).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A846").Value = 11
$ws.Range("B846").Value = "Vega Monumental Concepción"
$ws.Range("C846").Value = "Bíobío"
$ws.Range("D846").Value = 45267
$ws.Range("E846").Value = 8
$ws.Range("F846").Value = "Fruta"
$ws.Range("G846").Value = 100104
$ws.Range("H846").Value = "Frutos de pepita"
$ws.Range("I846").Value = 100104005
$ws.Range("J846").Value = "Pera"
$ws.Range("K846").Value = "Packham's Triumph"
$ws.Range("L846").Value = "Primera"
$ws.Range("M846").Value = 200
$ws.Range("N846").Value = 15000
$ws.Range("O846").Value = 15000
$ws.Range("P846").Value = 15000
$ws.Range("Q846").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R846").Value = "Región de O'Higgins"
$ws.Range("S846").Value = 833
$ws.Range("T846").Value = 18
